$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Cells.Item(33, 8).Value = 1504
$ws1.Cells.Item(33, 9).Value = 1768.75
$ws1.Cells.Item(33, 10).Value = 445
$ws1.Cells.Item(33, 11).Value = 1768.75
$ws1.Cells.Item(33, 12).Value = 445
$ws1.Cells.Item(33, 13).Value = -1539.75
$ws1.Cells.Item(33, 14).Value = -903

$ws1.Cells.Item(70, 8).Value = 1293.8536
$ws1.Cells.Item(70, 9).Value = 1308.5
$ws1.Cells.Item(70, 10).Value = 1268.4667
$ws1.Cells.Item(70, 11).Value = 3925.5
$ws1.Cells.Item(70, 12).Value = 3805.4001
$ws1.Cells.Item(70, 13).Value = -3655.5
$ws1.Cells.Item(70, 14).Value = -4345.4001

$ws1.Cells.Item(73, 8).Value = 1293.8536
$ws1.Cells.Item(73, 9).Value = 1308.5
$ws1.Cells.Item(73, 10).Value = 1268.4667
$ws1.Cells.Item(73, 11).Value = 3925.5
$ws1.Cells.Item(73, 12).Value = 3805.4001
$ws1.Cells.Item(73, 13).Value = -2989.5
$ws1.Cells.Item(73, 14).Value = -5677.4001

$ws1.Cells.Item(86, 8).Value = 137355.89
$ws1.Cells.Item(86, 9).Value = 205250.5
$ws1.Cells.Item(86, 10).Value = 1566.6666
$ws1.Cells.Item(86, 11).Value = 205250.5
$ws1.Cells.Item(86, 12).Value = 1566.6666
$ws1.Cells.Item(86, 13).Value = -204127.5
$ws1.Cells.Item(86, 14).Value = -3812.6666

$ws1.Cells.Item(88, 8).Value = 1674.8
$ws1.Cells.Item(88, 9).Value = 0
$ws1.Cells.Item(88, 10).Value = 1674.8
$ws1.Cells.Item(88, 11).Value = 0
$ws1.Cells.Item(88, 12).Value = 1674.8
$ws1.Cells.Item(88, 13).ClearContents()
$ws1.Cells.Item(88, 14).Value = -2486.8

$ws1.Cells.Item(89, 8).Value = 137355.89
$ws1.Cells.Item(89, 9).Value = 205250.5
$ws1.Cells.Item(89, 10).Value = 1566.6666
$ws1.Cells.Item(89, 11).Value = 1026252.5
$ws1.Cells.Item(89, 12).Value = 7833.333000000001
$ws1.Cells.Item(89, 13).Value = -1020636.5
$ws1.Cells.Item(89, 14).Value = -19065.333

$ws1.Cells.Item(91, 8).Value = 1674.8
$ws1.Cells.Item(91, 9).Value = 0
$ws1.Cells.Item(91, 10).Value = 1674.8
$ws1.Cells.Item(91, 11).Value = 0
$ws1.Cells.Item(91, 12).Value = 1674.8
$ws1.Cells.Item(91, 13).ClearContents()
$ws1.Cells.Item(91, 14).Value = -4482.8

$ws1.Cells.Item(105, 8).Value = 48373.332
$ws1.Cells.Item(105, 10).Value = 48373.332
$ws1.Cells.Item(105, 12).Value = 48373.332
$ws1.Cells.Item(105, 14).Value = -55361.332

$ws1.Cells.Item(132, 8).Value = 8072.1797
$ws1.Cells.Item(132, 9).Value = 9316.579
$ws1.Cells.Item(132, 10).Value = 6890
$ws1.Cells.Item(132, 11).Value = 27949.737
$ws1.Cells.Item(132, 12).Value = 20670
$ws1.Cells.Item(132, 13).Value = -25419.737
$ws1.Cells.Item(132, 14).Value = -25730

$ws1.Cells.Item(138, 8).Value = 2076.3914
$ws1.Cells.Item(138, 9).Value = 1831.7646
$ws1.Cells.Item(138, 10).Value = 2769.5
$ws1.Cells.Item(138, 11).Value = 5495.293799999999
$ws1.Cells.Item(138, 12).Value = 8308.5
$ws1.Cells.Item(138, 13).Value = -355.2937999999995
$ws1.Cells.Item(138, 14).Value = -18588.5

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Cells.Item(24, 8).Value = 24100
$ws2.Cells.Item(24, 10).Value = 24100
$ws2.Cells.Item(24, 12).Value = 24100
$ws2.Cells.Item(24, 14).Value = -24848

$ws2.Cells.Item(96, 8).Value = 26366.545
$ws2.Cells.Item(96, 10).Value = 26366.545
$ws2.Cells.Item(96, 12).Value = 26366.545
$ws2.Cells.Item(96, 14).Value = -31858.545

$ws2.Cells.Item(100, 8).Value = 24100
$ws2.Cells.Item(100, 10).Value = 24100
$ws2.Cells.Item(100, 12).Value = 24100
$ws2.Cells.Item(100, 14).Value = -26264

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Cells.Item(12, 8).Value = 552.5
$ws3.Cells.Item(12, 9).Value = 552.5
$ws3.Cells.Item(12, 11).Value = 552.5
$ws3.Cells.Item(12, 13).Value = -384.5

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Cells.Item(16, 8).Value = 3603
$ws4.Cells.Item(16, 9).Value = 3628.75
$ws4.Cells.Item(16, 10).Value = 3500
$ws4.Cells.Item(16, 11).Value = 3628.75
$ws4.Cells.Item(16, 12).Value = 3500
$ws4.Cells.Item(16, 13).Value = -3341.75
$ws4.Cells.Item(16, 14).Value = -4074

$ws4.Cells.Item(31, 8).Value = 2106.2642
$ws4.Cells.Item(31, 9).Value = 2007.4828
$ws4.Cells.Item(31, 10).Value = 2225.625
$ws4.Cells.Item(31, 11).Value = 2007.4828
$ws4.Cells.Item(31, 12).Value = 2225.625
$ws4.Cells.Item(31, 13).Value = -1712.4828
$ws4.Cells.Item(31, 14).Value = -2815.625

$ws4.Cells.Item(34, 8).Value = 2106.2642
$ws4.Cells.Item(34, 9).Value = 2007.4828
$ws4.Cells.Item(34, 10).Value = 2225.625
$ws4.Cells.Item(34, 11).Value = 2007.4828
$ws4.Cells.Item(34, 12).Value = 2225.625
$ws4.Cells.Item(34, 13).Value = -1805.4828
$ws4.Cells.Item(34, 14).Value = -2629.625

$ws4.Cells.Item(57, 8).Value = 8000
$ws4.Cells.Item(57, 10).Value = 8000
$ws4.Cells.Item(57, 12).Value = 8000
$ws4.Cells.Item(57, 14).Value = -9120

$ws4.Cells.Item(58, 8).Value = 4122.146
$ws4.Cells.Item(58, 9).Value = 3374.3845
$ws4.Cells.Item(58, 10).Value = 4399.8857
$ws4.Cells.Item(58, 11).Value = 3374.3845
$ws4.Cells.Item(58, 12).Value = 4399.8857
$ws4.Cells.Item(58, 13).Value = -3171.3845
$ws4.Cells.Item(58, 14).Value = -4805.8857

$ws4.Cells.Item(112, 8).Value = 50000
$ws4.Cells.Item(112, 10).Value = 50000
$ws4.Cells.Item(112, 12).Value = 50000
$ws4.Cells.Item(112, 14).Value = -52954

$ws4.Cells.Item(113, 8).Value = 3603
$ws4.Cells.Item(113, 9).Value = 3628.75
$ws4.Cells.Item(113, 10).Value = 3500
$ws4.Cells.Item(113, 11).Value = 3628.75
$ws4.Cells.Item(113, 12).Value = 3500
$ws4.Cells.Item(113, 13).Value = -1458.75
$ws4.Cells.Item(113, 14).Value = -7840

$ws4.Cells.Item(118, 8).Value = 78742
$ws4.Cells.Item(118, 10).Value = 78742
$ws4.Cells.Item(118, 12).Value = 78742
$ws4.Cells.Item(118, 14).Value = -82056

$ws4.Cells.Item(136, 8).Value = 4122.146
$ws4.Cells.Item(136, 9).Value = 3374.3845
$ws4.Cells.Item(136, 10).Value = 4399.8857
$ws4.Cells.Item(136, 11).Value = 10123.1535
$ws4.Cells.Item(136, 12).Value = 13199.6571
$ws4.Cells.Item(136, 13).Value = -7573.1535
$ws4.Cells.Item(136, 14).Value = -18299.6571

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Cells.Item(6, 8).Value = 1129
$ws5.Cells.Item(6, 9).Value = 161.25
$ws5.Cells.Item(6, 11).Value = 483.75
$ws5.Cells.Item(6, 13).Value = -370.75

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Cells.Item(113, 8).Value = 1180.0526
$ws6.Cells.Item(113, 9).Value = 1043.6666
$ws6.Cells.Item(113, 10).Value = 1413.8572
$ws6.Cells.Item(113, 11).Value = 1043.6666
$ws6.Cells.Item(113, 12).Value = 1413.8572
$ws6.Cells.Item(113, 13).Value = 1126.3334
$ws6.Cells.Item(113, 14).Value = -5753.8572

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Cells.Item(16, 8).Value = 2844.4443
$ws7.Cells.Item(16, 9).Value = 1200
$ws7.Cells.Item(16, 10).Value = 16000
$ws7.Cells.Item(16, 11).Value = 1200
$ws7.Cells.Item(16, 12).Value = 16000
$ws7.Cells.Item(16, 13).Value = -1030
$ws7.Cells.Item(16, 14).Value = -16340

$ws7.Cells.Item(46, 8).Value = 910289.0600000001
$ws7.Cells.Item(46, 9).Value = 800
$ws7.Cells.Item(46, 10).Value = 1251347.5
$ws7.Cells.Item(46, 11).Value = 800
$ws7.Cells.Item(46, 12).Value = 1251347.5
$ws7.Cells.Item(46, 13).Value = -612
$ws7.Cells.Item(46, 14).Value = -1251723.5

$ws7.Cells.Item(68, 8).Value = 2475.0833
$ws7.Cells.Item(68, 9).Value = 1965.125
$ws7.Cells.Item(68, 10).Value = 3495
$ws7.Cells.Item(68, 11).Value = 1965.125
$ws7.Cells.Item(68, 12).Value = 3495
$ws7.Cells.Item(68, 13).Value = -1216.125
$ws7.Cells.Item(68, 14).Value = -4993

$ws7.Cells.Item(71, 8).Value = 2475.0833
$ws7.Cells.Item(71, 9).Value = 1965.125
$ws7.Cells.Item(71, 10).Value = 3495
$ws7.Cells.Item(71, 11).Value = 9825.625
$ws7.Cells.Item(71, 12).Value = 17475
$ws7.Cells.Item(71, 13).Value = -6081.625
$ws7.Cells.Item(71, 14).Value = -24963

$ws7.Cells.Item(82, 8).Value = 1402.6522
$ws7.Cells.Item(82, 9).Value = 985.53845
$ws7.Cells.Item(82, 10).Value = 1944.9
$ws7.Cells.Item(82, 11).Value = 985.53845
$ws7.Cells.Item(82, 12).Value = 1944.9
$ws7.Cells.Item(82, 13).Value = -624.53845
$ws7.Cells.Item(82, 14).Value = -2666.9

$ws7.Cells.Item(85, 8).Value = 1402.6522
$ws7.Cells.Item(85, 9).Value = 985.53845
$ws7.Cells.Item(85, 10).Value = 1944.9
$ws7.Cells.Item(85, 11).Value = 985.53845
$ws7.Cells.Item(85, 12).Value = 1944.9
$ws7.Cells.Item(85, 13).Value = 262.46155
$ws7.Cells.Item(85, 14).Value = -4440.9

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Cells.Item(132, 8).Value = 1273.6394
$ws8.Cells.Item(132, 9).Value = 955.8222
$ws8.Cells.Item(132, 10).Value = 2167.5
$ws8.Cells.Item(132, 11).Value = 2867.4666
$ws8.Cells.Item(132, 12).Value = 6502.5
$ws8.Cells.Item(132, 13).Value = -337.4665999999997
$ws8.Cells.Item(132, 14).Value = -11562.5

$ws8.Cells.Item(136, 8).Value = 4103487.8
$ws8.Cells.Item(136, 9).Value = 6104223
$ws8.Cells.Item(136, 10).Value = 1980.5
$ws8.Cells.Item(136, 11).Value = 18312669
$ws8.Cells.Item(136, 12).Value = 5941.5
$ws8.Cells.Item(136, 13).Value = -18310119
$ws8.Cells.Item(136, 14).Value = -11041.5
